# Append new scrape results (2025-11-23 18:23 JST) to the top of the
# "ランサーズ" sheet, pushing the two existing rows down and interleaving
# three additional newly-scraped rows between/after them.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove the existing hyperlinks up front; we'll recreate them once all
# rows are in their final place so the relationship ids line up cleanly.
$ws.Hyperlinks.Delete()

# Insert three fresh rows. Doing this top-down means each insert shifts
# everything below (including rows inserted earlier in this loop) down
# by one, which reproduces the exact row layout from the diff:
#   old row2 (Java...)       -> row3
#   old row3 (WordPress...)  -> row5
$ws.Rows("2:2").Insert()
$ws.Rows("4:4").Insert()
$ws.Rows("6:6").Insert()

# --- Row 2 (brand new): data-scraping job -------------------------------
$ws.Range("A2").Value = "2025-11-23 18:23:41"
$ws.Range("B2").Value = "【急募】大規模データ収集自動化(スクレイピング・DB連携・エラー管理)案件"
$ws.Range("C2").Value = "システム開発"
$ws.Range("D2").Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Range("E2").Value = "期限情報なし"
$ws.Range("F2").Value = "https://www.lancers.jp/work/detail/5440052"
$ws.Range("G2").Value = 158
$ws.Range("H2").Value = "◆自動化,スクレイピング ◇管理"

# --- Row 3 (shifted from old row 2): only the timestamp changes --------
$ws.Range("A3").Value = "2025-11-23 18:23:41"

# --- Row 4 (brand new): Rakuten receipt tool job ------------------------
$ws.Range("A4").Value = "2025-11-23 18:23:41"
$ws.Range("B4").Value = "【急募】楽天市場の領収書一括ダウンロードツール開発依頼"
$ws.Range("C4").Value = "システム開発"
$ws.Range("D4").Value = "10,000 円 ~ 20,000 円 / 固定"
$ws.Range("E4").Value = "期限情報なし"
$ws.Range("F4").Value = "https://www.lancers.jp/work/detail/5440010"
$ws.Range("G4").Value = 120
$ws.Range("H4").Value = "◆ツール,開発"

# --- Row 5 (shifted from old row 3): only the timestamp changes --------
$ws.Range("A5").Value = "2025-11-23 18:23:41"

# --- Row 6 (brand new): program fix job (no skill summary) -------------
$ws.Range("A6").Value = "2025-11-23 18:23:41"
$ws.Range("B6").Value = "【急募】プログラム修正依頼!スキルを活かしてみませんか?"
$ws.Range("C6").Value = "システム開発"
$ws.Range("D6").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("E6").Value = "期限情報なし"
$ws.Range("F6").Value = "https://www.lancers.jp/work/detail/5440002"
$ws.Range("G6").Value = 13

# Recreate the URL hyperlinks (and their "Hyperlink" look) for every row.
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5440052")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5439921")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5440010")
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5439670")
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.lancers.jp/work/detail/5440002")
$ws.Range("F2:F6").Style = "Hyperlink"

# Widen the "スキル概要" (skill summary) column from 16 to 18 characters.
$ws.Columns.Item(8).ColumnWidth = 18 - 0.8333333333333334
